$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing name column values (string contents changed)
$ws.Range("B1").Value = "ahmed"
$ws.Range("B2").Value = "menna"

# Add gender column (P) and age column (Q)
$ws.Range("P1").Value = "male"
$ws.Range("P2").Value = "female"

$ws.Columns("Q").NumberFormat = "0"
$ws.Range("Q1").Value = 25
$ws.Range("Q2").Value = 22

# Match the width used by the other "bestFit" columns (M:O)
$ws.Columns("Q").ColumnWidth = $ws.Columns("M").ColumnWidth

# Update the selected cell to B3 (next empty row for the next patient)
$ws.Range("B3").Select() | Out-Null

# Set the sheet to print in portrait orientation
$ws.PageSetup.Orientation = 1
